$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H4").Value = 703.7143
$ws_ALC.Range("I4").Value = 654.3333
$ws_ALC.Range("K4").Value = 654.3333
$ws_ALC.Range("M4").Value = -540.3333
$ws_ALC.Range("H33").Value = 173.13333
$ws_ALC.Range("I33").Value = 153.84616
$ws_ALC.Range("J33").Value = 298.5
$ws_ALC.Range("K33").Value = 153.84616
$ws_ALC.Range("L33").Value = 298.5
$ws_ALC.Range("M33").Value = 75.15384
$ws_ALC.Range("N33").Value = -756.5
$ws_ALC.Range("H112").Value = 2083.4
$ws_ALC.Range("J112").Value = 2083.4
$ws_ALC.Range("L112").Value = 6250.200000000001
$ws_ALC.Range("N112").Value = -8466.200000000001
$ws_ALC.Range("H137").Value = 2493.1365
$ws_ALC.Range("I137").Value = 1829.9333
$ws_ALC.Range("J137").Value = 3914.2856
$ws_ALC.Range("K137").Value = 5489.7999
$ws_ALC.Range("L137").Value = 11742.8568
$ws_ALC.Range("M137").Value = -2939.7999
$ws_ALC.Range("N137").Value = -16842.8568
$ws_ALC.Range("H138").Value = 3889.8823
$ws_ALC.Range("J138").Value = 5423.8184
$ws_ALC.Range("L138").Value = 16271.4552
$ws_ALC.Range("N138").Value = -26551.4552

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H35").Value = 1720.2
$ws_ARM.Range("I35").Value = 1720.2
$ws_ARM.Range("K35").Value = 1720.2
$ws_ARM.Range("M35").Value = -1314.2
$ws_ARM.Range("H46").Value = 14992.083
$ws_ARM.Range("I46").Value = 15984.125
$ws_ARM.Range("J46").Value = 13008
$ws_ARM.Range("K46").Value = 15984.125
$ws_ARM.Range("L46").Value = 13008
$ws_ARM.Range("M46").Value = -15665.125
$ws_ARM.Range("N46").Value = -13646
$ws_ARM.Range("H61").Value = 2033.1666
$ws_ARM.Range("I61").Value = 1878.0625
$ws_ARM.Range("K61").Value = 1878.0625
$ws_ARM.Range("M61").Value = -1666.0625
$ws_ARM.Range("H74").Value = 1768.6522
$ws_ARM.Range("J74").Value = 4709.5
$ws_ARM.Range("L74").Value = 4709.5
$ws_ARM.Range("N74").Value = -6457.5
$ws_ARM.Range("H77").Value = 1768.6522
$ws_ARM.Range("J77").Value = 4709.5
$ws_ARM.Range("L77").Value = 23547.5
$ws_ARM.Range("N77").Value = -32283.5
$ws_ARM.Range("H122").Value = 402920.84
$ws_ARM.Range("I122").Value = 558393.2
$ws_ARM.Range("J122").Value = 3134.8572
$ws_ARM.Range("K122").Value = 1675179.6
$ws_ARM.Range("L122").Value = 9404.571599999999
$ws_ARM.Range("M122").Value = -1672729.6
$ws_ARM.Range("N122").Value = -14304.5716
$ws_ARM.Range("H132").Value = 2175
$ws_ARM.Range("J132").Value = 2249.25
$ws_ARM.Range("L132").Value = 6747.75
$ws_ARM.Range("N132").Value = -11807.75
$ws_ARM.Range("H136").Value = 2033.1666
$ws_ARM.Range("I136").Value = 1878.0625
$ws_ARM.Range("K136").Value = 5634.1875
$ws_ARM.Range("M136").Value = -3084.1875

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 658.03845
$ws_BSM.Range("J94").Value = 774.75
$ws_BSM.Range("L94").Value = 774.75
$ws_BSM.Range("N94").Value = -1676.75
$ws_BSM.Range("H105").Value = 2783.8235
$ws_BSM.Range("J105").Value = 5000
$ws_BSM.Range("L105").Value = 5000
$ws_BSM.Range("N105").Value = -8494
$ws_BSM.Range("H134").Value = 2530.45
$ws_BSM.Range("I134").Value = 2285.1667
$ws_BSM.Range("J134").Value = 2898.375
$ws_BSM.Range("K134").Value = 6855.500100000001
$ws_BSM.Range("L134").Value = 8695.125
$ws_BSM.Range("M134").Value = -4320.500100000001
$ws_BSM.Range("N134").Value = -13765.125

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 4526.8623
$ws_CRP.Range("I31").Value = 3555.158
$ws_CRP.Range("J31").Value = 6373.1
$ws_CRP.Range("K31").Value = 3555.158
$ws_CRP.Range("L31").Value = 6373.1
$ws_CRP.Range("M31").Value = -3260.158
$ws_CRP.Range("N31").Value = -6963.1
$ws_CRP.Range("H34").Value = 4526.8623
$ws_CRP.Range("I34").Value = 3555.158
$ws_CRP.Range("J34").Value = 6373.1
$ws_CRP.Range("K34").Value = 3555.158
$ws_CRP.Range("L34").Value = 6373.1
$ws_CRP.Range("M34").Value = -3353.158
$ws_CRP.Range("N34").Value = -6777.1
$ws_CRP.Range("H50").Value = 49999
$ws_CRP.Range("J50").Value = 49999
$ws_CRP.Range("L50").Value = 49999
$ws_CRP.Range("N50").Value = -51249
$ws_CRP.Range("H86").Value = 10397.111
$ws_CRP.Range("I86").Value = 8955
$ws_CRP.Range("J86").Value = 12199.75
$ws_CRP.Range("K86").Value = 8955
$ws_CRP.Range("L86").Value = 12199.75
$ws_CRP.Range("M86").Value = -7832
$ws_CRP.Range("N86").Value = -14445.75
$ws_CRP.Range("H89").Value = 10397.111
$ws_CRP.Range("I89").Value = 8955
$ws_CRP.Range("J89").Value = 12199.75
$ws_CRP.Range("K89").Value = 44775
$ws_CRP.Range("L89").Value = 60998.75
$ws_CRP.Range("M89").Value = -39159
$ws_CRP.Range("N89").Value = -72230.75
$ws_CRP.Range("H122").Value = 3009.6428
$ws_CRP.Range("I122").Value = 3010.4614
$ws_CRP.Range("K122").Value = 9031.3842
$ws_CRP.Range("M122").Value = -6581.3842
$ws_CRP.Range("H132").Value = 1900.5098
$ws_CRP.Range("I132").Value = 1581.875
$ws_CRP.Range("K132").Value = 4745.625
$ws_CRP.Range("M132").Value = -2215.625
$ws_CRP.Range("H134").Value = 2683.25
$ws_CRP.Range("I134").Value = 2197.2727
$ws_CRP.Range("K134").Value = 6591.8181
$ws_CRP.Range("M134").Value = -4056.8181

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H38").Value = 119.875
$ws_CUL.Range("J38").Value = 191.25
$ws_CUL.Range("L38").Value = 573.75
$ws_CUL.Range("N38").Value = -1267.75
$ws_CUL.Range("H132").Value = 11165.25
$ws_CUL.Range("J132").Value = 11333
$ws_CUL.Range("L132").Value = 101997
$ws_CUL.Range("N132").Value = -107057

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 2187.8918
$ws_GSM.Range("J132").Value = 2449.476
$ws_GSM.Range("L132").Value = 7348.428
$ws_GSM.Range("N132").Value = -12408.428

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 2943.5557
$ws_LTW.Range("I40").Value = 2943.5557
$ws_LTW.Range("K40").Value = 2943.5557
$ws_LTW.Range("M40").Value = -2807.5557
$ws_LTW.Range("H46").Value = 3498.625
$ws_LTW.Range("I46").Value = 3331.6667
$ws_LTW.Range("J46").Value = 3999.5
$ws_LTW.Range("K46").Value = 3331.6667
$ws_LTW.Range("L46").Value = 3999.5
$ws_LTW.Range("M46").Value = -3143.6667
$ws_LTW.Range("N46").Value = -4375.5
$ws_LTW.Range("H61").Value = 5499.6
$ws_LTW.Range("I61").Value = 5999.5
$ws_LTW.Range("K61").Value = 5999.5
$ws_LTW.Range("M61").Value = -5797.5
$ws_LTW.Range("H68").Value = 1335.4286
$ws_LTW.Range("I68").Value = 1349.8
$ws_LTW.Range("J68").Value = 1299.5
$ws_LTW.Range("K68").Value = 1349.8
$ws_LTW.Range("L68").Value = 1299.5
$ws_LTW.Range("M68").Value = -600.8
$ws_LTW.Range("N68").Value = -2797.5
$ws_LTW.Range("H71").Value = 1335.4286
$ws_LTW.Range("I71").Value = 1349.8
$ws_LTW.Range("J71").Value = 1299.5
$ws_LTW.Range("K71").Value = 6749
$ws_LTW.Range("L71").Value = 6497.5
$ws_LTW.Range("M71").Value = -3005
$ws_LTW.Range("N71").Value = -13985.5
$ws_LTW.Range("H76").Value = 32666.666
$ws_LTW.Range("J76").Value = 32666.666
$ws_LTW.Range("L76").Value = 32666.666
$ws_LTW.Range("N76").Value = -33342.666
$ws_LTW.Range("H79").Value = 32666.666
$ws_LTW.Range("J79").Value = 32666.666
$ws_LTW.Range("L79").Value = 32666.666
$ws_LTW.Range("N79").Value = -35006.666
$ws_LTW.Range("H82").Value = 1951.7858
$ws_LTW.Range("I82").Value = 2092.6
$ws_LTW.Range("K82").Value = 2092.6
$ws_LTW.Range("M82").Value = -1731.6
$ws_LTW.Range("H85").Value = 1951.7858
$ws_LTW.Range("I85").Value = 2092.6
$ws_LTW.Range("K85").Value = 2092.6
$ws_LTW.Range("M85").Value = -844.5999999999999
$ws_LTW.Range("H101").Value = 38000
$ws_LTW.Range("J101").Value = 38000
$ws_LTW.Range("L101").Value = 38000
$ws_LTW.Range("N101").Value = -44490
$ws_LTW.Range("H113").Value = 5499.6
$ws_LTW.Range("I113").Value = 5999.5
$ws_LTW.Range("K113").Value = 5999.5
$ws_LTW.Range("M113").Value = -3829.5
$ws_LTW.Range("I122").Value = 9330.333000000001
$ws_LTW.Range("J122").Value = 5091.75
$ws_LTW.Range("K122").Value = 27990.999
$ws_LTW.Range("L122").Value = 15275.25
$ws_LTW.Range("M122").Value = -25540.999
$ws_LTW.Range("N122").Value = -20175.25
$ws_LTW.Range("H132").Value = 3853.4644
$ws_LTW.Range("I132").Value = 2919.85
$ws_LTW.Range("K132").Value = 8759.549999999999
$ws_LTW.Range("M132").Value = -6229.549999999999

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H45").Value = 13153.375
$ws_WVR.Range("I45").Value = 14991
$ws_WVR.Range("J45").Value = 12890.857
$ws_WVR.Range("K45").Value = 14991
$ws_WVR.Range("L45").Value = 12890.857
$ws_WVR.Range("M45").Value = -14500
$ws_WVR.Range("N45").Value = -13872.857
$ws_WVR.Range("H62").Value = 7362.25
$ws_WVR.Range("I62").Value = 6374.75
$ws_WVR.Range("J62").Value = 8349.75
$ws_WVR.Range("K62").Value = 6374.75
$ws_WVR.Range("L62").Value = 8349.75
$ws_WVR.Range("M62").Value = -5750.75
$ws_WVR.Range("N62").Value = -9597.75
$ws_WVR.Range("H65").Value = 7362.25
$ws_WVR.Range("I65").Value = 6374.75
$ws_WVR.Range("J65").Value = 8349.75
$ws_WVR.Range("K65").Value = 31873.75
$ws_WVR.Range("L65").Value = 41748.75
$ws_WVR.Range("M65").Value = -28753.75
$ws_WVR.Range("N65").Value = -47988.75
$ws_WVR.Range("H107").Value = 663.55554
$ws_WVR.Range("I107").Value = 710.2857
$ws_WVR.Range("K107").Value = 2130.8571
$ws_WVR.Range("M107").Value = -210.8571000000002
$ws_WVR.Range("H112").Value = 17187
$ws_WVR.Range("J112").Value = 17187
$ws_WVR.Range("L112").Value = 17187
$ws_WVR.Range("N112").Value = -20141
$ws_WVR.Range("H122").Value = 3880.375
$ws_WVR.Range("I122").Value = 4673.8335
$ws_WVR.Range("K122").Value = 14021.5005
$ws_WVR.Range("M122").Value = -11571.5005
$ws_WVR.Range("H132").Value = 48498
$ws_WVR.Range("I132").Value = 60141.812
$ws_WVR.Range("K132").Value = 180425.436
$ws_WVR.Range("M132").Value = -177895.436

Write-Host "Applied 235 cell updates across 8 sheets"
